# Apply the diff:
#  1. Clear the stray empty inline-string placeholder cells on Table_1
#     (B2, A3, B37) so they no longer emit a <c> element.
#  2. Add a new worksheet "Table_2" right after "Table_1" holding the
#     capital-adequacy ratio table (header row + 3 data rows).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Table_1")

# --- 1. Remove the empty placeholder cells -------------------------------
$ws1.Range("B2").Clear()
$ws1.Range("A3").Clear()
$ws1.Range("B37").Clear()

# --- 2. Add "Table_2" positioned after "Table_1" --------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Table_2"

# Header row text.
$ws2.Range("A1").Value = "Əmsal"
$ws2.Range("B1").Value = "Norma (Sistem əhəmiyyətli)"
$ws2.Range("C1").Value = "Norma (Banklar istisna)"
$ws2.Range("D1").Value = "Fakt"

# Match Table_1's header formatting (bold, bordered, centered/top-aligned)
# by copying the existing header style instead of re-building it, so the
# workbook doesn't pick up a near-duplicate font/xf.
$ws1.Range("A1:B1").Copy()
$ws2.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows. Force text format first so the percentage-looking strings are
# kept as literal text ("6.0%") instead of being auto-converted to numeric
# percentages by Excel's data recognition.
$dataRange = $ws2.Range("A2:D4")
$dataRange.NumberFormat = "@"

$ws2.Range("A2").Value = "9.  I dərəcəli  kapitalın  adekvatlıq əmsalı"
$ws2.Range("B2").Value = "6.0%"
$ws2.Range("C2").Value = "5.0%"
$ws2.Range("D2").Value = "10.0%"

$ws2.Range("A3").Value = "10. məcmu kapitalın  adekvatlıq  əmsalı"
$ws2.Range("B3").Value = "12.0%"
$ws2.Range("C3").Value = "10.0%"
$ws2.Range("D3").Value = "15.0%"

$ws2.Range("A4").Value = "11. Leverec əmsalı"
$ws2.Range("B4").Value = "minimum 5%"
$ws2.Range("C4").Value = "minimum 4%"
$ws2.Range("D4").Value = "6.0%"

# Now that the text is locked in, restore the cells' formatting back to the
# workbook's plain default style (an untouched cell carries that style) so
# the new rows don't end up with a stray "Text" number format applied.
$ws1.Range("C1").Copy()
$dataRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Keep Table_1 as the active sheet/tab (matches the original book view).
$ws1.Activate()
